# Apply "Code updated with cookie handling" changes to the instamart product sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 11: Too Yumm! B2G1 Potato Chips - Indian Masala (Pursaiwakam, 600007/600008)
# MRP/Price swap and variants collapsed to "No variants"
# ---------------------------------------------------------------------------
$ws.Range("C11").Value = 70
$ws.Range("D11").Value = 49
$ws.Range("F11").Value = "[['No variants']]"

# ---------------------------------------------------------------------------
# Row 13: Go Zero Belgian Dark Chocolate Cup (Pursaiwakam) -- C/D/H become numeric
# ---------------------------------------------------------------------------
$ws.Range("C13").NumberFormat = "General"
$ws.Range("C13").Value = 120
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Value = 98
$ws.Range("H13").NumberFormat = "General"
$ws.Range("H13").Value = 600008

# ---------------------------------------------------------------------------
# New rows 14-22
# ---------------------------------------------------------------------------

function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# Row 14: Details not available (New Karol Bagh store, pincode 110006)
Set-TextCell $ws.Range("A14") "Details not available for this location"
Set-TextCell $ws.Range("B14") "Details not available for this location"
Set-TextCell $ws.Range("C14") "Details not available for this location"
Set-TextCell $ws.Range("D14") "Details not available for this location"
Set-TextCell $ws.Range("E14") "Details not available for this location"
Set-TextCell $ws.Range("F14") "['Details not available for this location']"
Set-TextCell $ws.Range("G14") "Details not available for this location"
$ws.Range("H14").NumberFormat = "General"
$ws.Range("H14").Value = 110006

# Row 15: Go Zero Madagascar Chocobar (New Karol Bagh, 110006)
Set-TextCell $ws.Range("A15") "Go Zero Madagascar Chocobar Low Calorie Guilt Free Ice Cream Bar"
Set-TextCell $ws.Range("B15") "Go Zero"
Set-TextCell $ws.Range("C15") "135"
Set-TextCell $ws.Range("D15") "103"
Set-TextCell $ws.Range("E15") "Seller Name: PYD Retail Pvt Ltd - New Karol Bagh FSSAI Number: 20250325107131904 Address: UG Flr ,Plot 10130, Katra chajju pandit, main Rani Jhasi Road, New Delhi-110005"
Set-TextCell $ws.Range("F15") "[[{'Variant Size': '60 ml x 4', 'Price': '400', 'MRP': '500', 'Discount': '20% OFF', 'Status': 'In stock'}, {'Variant Size': '60 ml x 2', 'Price': '202', 'MRP': '250', 'Discount': '19% OFF', 'Status': 'In stock'}, {'Variant Size': '60 ml', 'Price': '102', 'MRP': '125', 'Discount': '18% OFF', 'Status': 'In stock'}]]"
Set-TextCell $ws.Range("G15") "https://www.swiggy.com/instamart/item/LNSIN80D9X?storeId=1402050"
$ws.Range("H15").NumberFormat = "General"
$ws.Range("H15").Value = 110006

# Row 16: Go Zero Belgian Dark Chocolate Cup (New Karol Bagh, 110006)
Set-TextCell $ws.Range("A16") "Go Zero Belgian Dark Chocolate Low Calorie Guilt Free Ice Cream Cup"
Set-TextCell $ws.Range("B16") "Go Zero"
Set-TextCell $ws.Range("C16") "120"
Set-TextCell $ws.Range("D16") "98"
Set-TextCell $ws.Range("E16") "Seller Name: PYD Retail Pvt Ltd - New Karol Bagh FSSAI Number: 20250325107131904 Address: UG Flr ,Plot 10130, Katra chajju pandit, main Rani Jhasi Road, New Delhi-110005"
Set-TextCell $ws.Range("F16") "[[{'Variant Size': '500 ml', 'Price': '327', 'MRP': '425', 'Discount': '23% OFF', 'Status': 'In stock'}, {'Variant Size': '100 ml x 4', 'Price': '379', 'MRP': '480', 'Discount': '21% OFF', 'Status': 'In stock'}, {'Variant Size': '100 ml', 'Price': '98', 'MRP': '120', 'Discount': '18% OFF', 'Status': 'In stock'}]]"
Set-TextCell $ws.Range("G16") "https://www.swiggy.com/instamart/item/YRL5V0ED04?storeId=1402050"
$ws.Range("H16").NumberFormat = "General"
$ws.Range("H16").Value = 110006

# Row 17: Details not available (Patel Nagar store, pincode 110008)
Set-TextCell $ws.Range("A17") "Details not available for this location"
Set-TextCell $ws.Range("B17") "Details not available for this location"
Set-TextCell $ws.Range("C17") "Details not available for this location"
Set-TextCell $ws.Range("D17") "Details not available for this location"
Set-TextCell $ws.Range("E17") "Details not available for this location"
Set-TextCell $ws.Range("F17") "['Details not available for this location']"
Set-TextCell $ws.Range("G17") "Details not available for this location"
$ws.Range("H17").NumberFormat = "General"
$ws.Range("H17").Value = 110008

# Row 18: Go Zero Madagascar Chocobar (Patel Nagar, 110008)
Set-TextCell $ws.Range("A18") "Go Zero Madagascar Chocobar Low Calorie Guilt Free Ice Cream Bar"
Set-TextCell $ws.Range("B18") "Go Zero"
Set-TextCell $ws.Range("C18") "135"
Set-TextCell $ws.Range("D18") "103"
Set-TextCell $ws.Range("E18") "Seller Name: PYD Retail Pvt. Ltd. Patel Nagar FSSAI Number: 13322006000126 Address: Plot No BP-08,West Patel Nagar,New Delhi-110008"
Set-TextCell $ws.Range("F18") "[[{'Variant Size': '60 ml x 4', 'Price': '400', 'MRP': '500', 'Discount': '20% OFF', 'Status': 'In stock'}, {'Variant Size': '60 ml x 2', 'Price': '202', 'MRP': '250', 'Discount': '19% OFF', 'Status': 'In stock'}, {'Variant Size': '60 ml', 'Price': '102', 'MRP': '125', 'Discount': '18% OFF', 'Status': 'In stock'}]]"
Set-TextCell $ws.Range("G18") "https://www.swiggy.com/instamart/item/LNSIN80D9X?storeId=1402050"
$ws.Range("H18").NumberFormat = "General"
$ws.Range("H18").Value = 110008

# Row 19: Go Zero Belgian Dark Chocolate Cup (Patel Nagar, 110008)
Set-TextCell $ws.Range("A19") "Go Zero Belgian Dark Chocolate Low Calorie Guilt Free Ice Cream Cup"
Set-TextCell $ws.Range("B19") "Go Zero"
Set-TextCell $ws.Range("C19") "120"
Set-TextCell $ws.Range("D19") "98"
Set-TextCell $ws.Range("E19") "Seller Name: PYD Retail Pvt. Ltd. Patel Nagar FSSAI Number: 13322006000126 Address: Plot No BP-08,West Patel Nagar,New Delhi-110008"
Set-TextCell $ws.Range("F19") "[[{'Variant Size': '500 ml', 'Price': '327', 'MRP': '425', 'Discount': '23% OFF', 'Status': 'In stock'}, {'Variant Size': '100 ml x 4', 'Price': '379', 'MRP': '480', 'Discount': '21% OFF', 'Status': 'In stock'}, {'Variant Size': '100 ml', 'Price': '98', 'MRP': '120', 'Discount': '18% OFF', 'Status': 'In stock'}]]"
Set-TextCell $ws.Range("G19") "https://www.swiggy.com/instamart/item/YRL5V0ED04?storeId=1402050"
$ws.Range("H19").NumberFormat = "General"
$ws.Range("H19").Value = 110008

# Row 20: Too Yumm! Chips Indian Masala (Jharsa Village, 122003)
Set-TextCell $ws.Range("A20") "Too Yumm! Chips Indian Masala"
Set-TextCell $ws.Range("B20") "Too Yumm!"
Set-TextCell $ws.Range("C20") "70"
Set-TextCell $ws.Range("D20") "49"
Set-TextCell $ws.Range("E20") "Seller Name: PYD Retail Pvt Ltd - Jharsa Village FSSAI Number: 20250116106886225 Address: GF, B-12, sector-45, near greenwood city, Gurgaon-122003"
Set-TextCell $ws.Range("F20") "[['No variants']]"
Set-TextCell $ws.Range("G20") "https://www.swiggy.com/instamart/item/N5SBE9SBEE?storeId=1402050"
$ws.Range("H20").NumberFormat = "General"
$ws.Range("H20").Value = 122003

# Row 21: Go Zero Madagascar Chocobar (Sector 50, 122003)
Set-TextCell $ws.Range("A21") "Go Zero Madagascar Chocobar Low Calorie Guilt Free Ice Cream Bar"
Set-TextCell $ws.Range("B21") "Go Zero"
Set-TextCell $ws.Range("C21") "135"
Set-TextCell $ws.Range("D21") "109"
Set-TextCell $ws.Range("E21") "Seller Name: PYD Retail Pvt Ltd - Sector 50 FSSAI Number: 20240903106414148 Address: U.no 105,106,108, 04, 110-112, gf to 7th floor Revenue Estate of Village Badshahpur, Sector 50, Gurgaon, Haryana, 122108"
Set-TextCell $ws.Range("F21") "[[{'Variant Size': '60 ml x 2', 'Price': '202', 'MRP': '250', 'Discount': '19% OFF', 'Status': 'In stock'}, {'Variant Size': '60 ml', 'Price': '102', 'MRP': '125', 'Discount': '18% OFF', 'Status': 'In stock'}]]"
Set-TextCell $ws.Range("G21") "https://www.swiggy.com/instamart/item/LNSIN80D9X?storeId=1402050"
$ws.Range("H21").NumberFormat = "General"
$ws.Range("H21").Value = 122003

# Row 22: Go Zero Belgian Dark Chocolate Cup (Sector 50, 122003) -- H22 stays textual
Set-TextCell $ws.Range("A22") "Go Zero Belgian Dark Chocolate Low Calorie Guilt Free Ice Cream Cup"
Set-TextCell $ws.Range("B22") "Go Zero"
Set-TextCell $ws.Range("C22") "425"
Set-TextCell $ws.Range("D22") "327"
Set-TextCell $ws.Range("E22") "Seller Name: PYD Retail Pvt Ltd - Sector 50 FSSAI Number: 20240903106414148 Address: U.no 105,106,108, 04, 110-112, gf to 7th floor Revenue Estate of Village Badshahpur, Sector 50, Gurgaon, Haryana, 122108"
Set-TextCell $ws.Range("F22") "[[{'Variant Size': '500 ml', 'Price': '327', 'MRP': '425', 'Discount': '23% OFF', 'Status': 'In stock'}, {'Variant Size': '100 ml x 4', 'Price': '379', 'MRP': '480', 'Discount': '21% OFF', 'Status': 'In stock'}, {'Variant Size': '100 ml', 'Price': '98', 'MRP': '120', 'Discount': '18% OFF', 'Status': 'Sold Out'}]]"
Set-TextCell $ws.Range("G22") "https://www.swiggy.com/instamart/item/YRL5V0ED04?storeId=1402050"
Set-TextCell $ws.Range("H22") "122003"
